$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

$ws.Range("AP1").Value = "VL"
$ws.Range("AQ1").Value = "SL"

$ws.Range("AP2:AP6").Value = "15"
$ws.Range("AQ2:AQ6").Value = "15"

$ws.Range("AP1:AQ6").Select() | Out-Null
